$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.474368810653687
$ws.Range("B1").Value = 6.31640100479126
$ws.Range("C1").Value = 5.840932846069336
$ws.Range("D1").Value = 6.688359260559082
$ws.Range("E1").Value = 3.737930536270142
